$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new row 11: label "F1_Score" in A11, and F1 formula in B11
$ws.Range("A11").Value = "F1_Score"
$ws.Range("B11").Formula = "=(2*B9*B10)/(B9+B10)"

# Match number format of B9/B10 (percentage style) for the new B11 cell
$ws.Range("B11").NumberFormat = $ws.Range("B9").NumberFormat

# Update the active selection to match the target state
$ws.Range("C17").Select()
